$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows right before the current row 8, shifting existing
# rows 8-25 down to 10-27 (matches dimension growing from A1:R25 to A1:R27).
$ws.Range("A8:A9").EntireRow.Insert()

# New row 8: Perejil, Primera, fecha 44859 (2022-10-25)
$ws.Cells.Item(8, 1).Value = 7
$ws.Cells.Item(8, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(8, 3).Value = "Ñuble"
$ws.Cells.Item(8, 4).Value = 44859
$ws.Cells.Item(8, 5).Value = 16
$ws.Cells.Item(8, 6).Value = 100112044
$ws.Cells.Item(8, 7).Value = "Perejil"
$ws.Cells.Item(8, 8).Value = "Sin especificar"
$ws.Cells.Item(8, 9).Value = "Primera"
$ws.Cells.Item(8, 10).Value = 300
$ws.Cells.Item(8, 11).Value = 700
$ws.Cells.Item(8, 12).Value = 800
$ws.Cells.Item(8, 13).Value = 750
$ws.Cells.Item(8, 14).Value = "$/atado 0,5 a 1 kilo"
$ws.Cells.Item(8, 15).Value = "Región del Maule"
$ws.Cells.Item(8, 16).Value = 750
$ws.Cells.Item(8, 17).Value = 1
$ws.Cells.Item(8, 18).Value = "Hortaliza"

# New row 9: Perejil, Segunda, fecha 44859 (2022-10-25)
$ws.Cells.Item(9, 1).Value = 7
$ws.Cells.Item(9, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(9, 3).Value = "Ñuble"
$ws.Cells.Item(9, 4).Value = 44859
$ws.Cells.Item(9, 5).Value = 16
$ws.Cells.Item(9, 6).Value = 100112044
$ws.Cells.Item(9, 7).Value = "Perejil"
$ws.Cells.Item(9, 8).Value = "Sin especificar"
$ws.Cells.Item(9, 9).Value = "Segunda"
$ws.Cells.Item(9, 10).Value = 200
$ws.Cells.Item(9, 11).Value = 600
$ws.Cells.Item(9, 12).Value = 600
$ws.Cells.Item(9, 13).Value = 600
$ws.Cells.Item(9, 14).Value = "$/atado 0,5 a 1 kilo"
$ws.Cells.Item(9, 15).Value = "Región del Maule"
$ws.Cells.Item(9, 16).Value = 600
$ws.Cells.Item(9, 17).Value = 1
$ws.Cells.Item(9, 18).Value = "Hortaliza"
